$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.656.03'
$ws.Range("E2").Value = '  -1.64%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.588.82'
$ws.Range("E3").Value = '  -2.17%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '210.76'
$ws.Range("E5").Value = '  -1.75%  '
$ws.Range("E6").Value = '  -2.10%  '
$ws.Range("E7").Value = '  +0.14%  '
$ws.Range("E8").Value = '  -1.96%  '
$ws.Range("E9").Value = '  -1.90%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.61'
$ws.Range("E10").Value = '  -3.45%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0832'
$ws.Range("E11").Value = '  -1.64%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.811.24'
$ws.Range("E12").Value = '  -2.16%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.587.89'
$ws.Range("E13").Value = '  -2.26%  '
$ws.Range("E14").Value = '  -2.64%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.521'
$ws.Range("E15").Value = '  -3.92%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '64.79'
$ws.Range("E16").Value = '  +0.41%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '26.624.11'
$ws.Range("E17").Value = '  -1.60%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.0₃0726'
$ws.Range("E18").Value = '  -2.08%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '208.50'
$ws.Range("E19").Value = '  -3.61%  '
$ws.Range("E20").Value = '  +0.18%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.72'
$ws.Range("E21").Value = '  -2.79%  '
$ws.Range("E22").Value = '  -3.20%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.34'
$ws.Range("E23").Value = '  -3.10%  '
$ws.Range("E24").Value = '  -1.94%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '146.86'
$ws.Range("E25").Value = '  -0.25%  '
$ws.Range("E26").Value = '  +0.18%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.23'
$ws.Range("E27").Value = '  -0.98%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.114'
$ws.Range("E28").Value = '  -3.37%  '
$ws.Range("E29").Value = '  -1.87%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0506'
$ws.Range("E30").Value = '  +0.11%  '
$ws.Range("E31").Value = '  -1.89%  '
$ws.Range("E32").Value = '  -3.96%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.665'
$ws.Range("E33").Value = '  +21.58%  '
$ws.Range("E34").Value = '  -3.00%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.306.61'
$ws.Range("E35").Value = '  -2.49%  '
$ws.Range("E36").Value = '  -1.25%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.49'
$ws.Range("E37").Value = '  -5.23%  '
$ws.Range("E38").Value = '  -2.91%  '
$ws.Range("E39").Value = '  -2.51%  '
$ws.Range("E40").Value = '  +0.16%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.793'
$ws.Range("E41").Value = '  -1.26%  '
$ws.Range("E42").Value = '  +2.34%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.17'
$ws.Range("E43").Value = '  -2.85%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '62.64'
$ws.Range("E44").Value = '  -4.26%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.724.82'
$ws.Range("E45").Value = '  -1.95%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '89.71'
$ws.Range("E46").Value = '  -0.87%  '
$ws.Range("E47").Value = '  -0.36%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.832'
$ws.Range("E48").Value = '  -2.06%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0980'
$ws.Range("E49").Value = '  -1.79%  '
$ws.Range("E50").Value = '  -1.79%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.54'
$ws.Range("E51").Value = '  -0.31%  '
